# Auto-generated edit script applying the Marilith_Profits.xlsx cell-value diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2635
$ws.Range("I2").Value = 1012.75
$ws.Range("J2").Value = 4798
$ws.Range("K2").Value = 1012.75
$ws.Range("L2").Value = 4798
$ws.Range("M2").Value = -899.75
$ws.Range("N2").Value = -5024

$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 20000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 20000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -21944

$ws.Range("H74").Value = 169333.33
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 500000
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 500000
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -501872

$ws.Range("H77").Value = 169333.33
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 500000
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 2500000
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -2509360

$ws.Range("H92").Value = 1355.625
$ws.Range("J92").Value = 112.5
$ws.Range("L92").Value = 112.5
$ws.Range("N92").Value = -2608.5

$ws.Range("H96").Value = 669.8570999999999
$ws.Range("I96").Value = 623.75
$ws.Range("J96").Value = 731.3333
$ws.Range("K96").Value = 1871.25
$ws.Range("L96").Value = 2193.9999
$ws.Range("M96").Value = -498.25
$ws.Range("N96").Value = -4939.9999

$ws.Range("H100").Value = 2334.9285
$ws.Range("I100").Value = 587.625
$ws.Range("J100").Value = 4664.6665
$ws.Range("K100").Value = 587.625
$ws.Range("L100").Value = 4664.6665
$ws.Range("M100").Value = -46.625
$ws.Range("N100").Value = -5746.6665

$ws.Range("H101").Value = 509.4
$ws.Range("I101").Value = 501.125
$ws.Range("J101").Value = 542.5
$ws.Range("K101").Value = 1503.375
$ws.Range("L101").Value = 1627.5
$ws.Range("M101").Value = 118.625
$ws.Range("N101").Value = -4871.5

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H129").Value = 1835.8889
$ws.Range("J129").Value = 3497.6365
$ws.Range("L129").Value = 10492.9095
$ws.Range("N129").Value = -20492.9095

$ws.Range("H138").Value = 1523.9166
$ws.Range("I138").Value = 571.75
$ws.Range("K138").Value = 1715.25
$ws.Range("M138").Value = 3424.75


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1108.5714
$ws.Range("I97").Value = 911
$ws.Range("K97").Value = 911
$ws.Range("M97").Value = -415

$ws.Range("H105").Value = 20750
$ws.Range("J105").Value = 20750
$ws.Range("L105").Value = 20750
$ws.Range("N105").Value = -27738

$ws.Range("H132").Value = 913
$ws.Range("I132").Value = 913
$ws.Range("K132").Value = 2739
$ws.Range("M132").Value = -209


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2797.8333
$ws.Range("I86").Value = 2797.8333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2797.8333
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1674.8333
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 2797.8333
$ws.Range("I89").Value = 2797.8333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 13989.1665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -8373.166499999999
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 2172.0715
$ws.Range("I94").Value = 2410
$ws.Range("J94").Value = 1299.6666
$ws.Range("K94").Value = 2410
$ws.Range("L94").Value = 1299.6666
$ws.Range("M94").Value = -1959
$ws.Range("N94").Value = -2201.6666

$ws.Range("H99").Value = 999
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 999
$ws.Range("M99").Value = 499

$ws.Range("H134").Value = 3100.9312
$ws.Range("I134").Value = 3100.9312
$ws.Range("K134").Value = 9302.793600000001
$ws.Range("M134").Value = -6767.793600000001


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4566.3335
$ws.Range("I31").Value = 4566.3335
$ws.Range("K31").Value = 4566.3335
$ws.Range("M31").Value = -4271.3335

$ws.Range("H34").Value = 4566.3335
$ws.Range("I34").Value = 4566.3335
$ws.Range("K34").Value = 4566.3335
$ws.Range("M34").Value = -4364.3335

$ws.Range("H99").Value = 4797.8
$ws.Range("I99").Value = 4499.5
$ws.Range("K99").Value = 4499.5
$ws.Range("M99").Value = -3001.5

$ws.Range("H122").Value = 1312.5
$ws.Range("I122").Value = 1080
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 3240
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -790
$ws.Range("N122").Value = -10000

$ws.Range("H126").Value = 4797.8
$ws.Range("I126").Value = 4499.5
$ws.Range("K126").Value = 13498.5
$ws.Range("M126").Value = -11028.5


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 196.1
$ws.Range("J12").Value = 298.33334
$ws.Range("L12").Value = 895.0000200000001
$ws.Range("N12").Value = -1241.00002

$ws.Range("H70").Value = 1500
$ws.Range("J70").Value = 1000
$ws.Range("L70").Value = 3000
$ws.Range("N70").Value = -3630

$ws.Range("H73").Value = 1500
$ws.Range("J73").Value = 1000
$ws.Range("L73").Value = 3000
$ws.Range("N73").Value = -5184

$ws.Range("H75").Value = 5548.75
$ws.Range("I75").Value = 2200
$ws.Range("J75").Value = 6665
$ws.Range("K75").Value = 6600
$ws.Range("L75").Value = 19995
$ws.Range("M75").Value = -5602
$ws.Range("N75").Value = -21991

$ws.Range("H78").Value = 5548.75
$ws.Range("I78").Value = 2200
$ws.Range("J78").Value = 6665
$ws.Range("K78").Value = 19800
$ws.Range("L78").Value = 59985
$ws.Range("M78").Value = -14808
$ws.Range("N78").Value = -69969


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 788.4286
$ws.Range("I97").Value = 622.7
$ws.Range("K97").Value = 622.7
$ws.Range("M97").Value = -126.7

$ws.Range("H113").Value = 2837.5
$ws.Range("I113").Value = 2837.5
$ws.Range("K113").Value = 2837.5
$ws.Range("M113").Value = -667.5

$ws.Range("H122").Value = 20836744
$ws.Range("I122").Value = 20836744
$ws.Range("K122").Value = 62510232
$ws.Range("M122").Value = -62507782

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5540.6665
$ws.Range("I40").Value = 3311
$ws.Range("K40").Value = 3311
$ws.Range("M40").Value = -3175

$ws.Range("H46").Value = 1559.8
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1559.8
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1559.8
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1935.8

$ws.Range("H61").Value = 7449.5
$ws.Range("I61").Value = 7449.5
$ws.Range("K61").Value = 7449.5
$ws.Range("M61").Value = -7247.5

$ws.Range("H93").Value = 1149.75
$ws.Range("I93").Value = 866.3333
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 866.3333
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 381.6667
$ws.Range("N93").Value = -4496

$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959

$ws.Range("H113").Value = 7449.5
$ws.Range("I113").Value = 7449.5
$ws.Range("K113").Value = 7449.5
$ws.Range("M113").Value = -5279.5


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1192.5
$ws.Range("I96").Value = 1192.5
$ws.Range("K96").Value = 1192.5
$ws.Range("M96").Value = 180.5

$ws.Range("H100").Value = 990
$ws.Range("I100").Value = 990
$ws.Range("K100").Value = 1980
$ws.Range("M100").Value = -1439

$ws.Range("H107").Value = 492.66666
$ws.Range("I107").Value = 529.25
$ws.Range("K107").Value = 1587.75
$ws.Range("M107").Value = 332.25

$ws.Range("H113").Value = 1470
$ws.Range("J113").Value = 1475
$ws.Range("L113").Value = 4425
$ws.Range("N113").Value = -8765

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H136").Value = 1093.0667
$ws.Range("I136").Value = 1093.0667
$ws.Range("K136").Value = 3279.2001
$ws.Range("M136").Value = -729.2001

